# Add a new "break_on_off" column (L) to the schedule sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header for new column L, matching style of existing header row (bold/centered = style used by A1:K1)
$ws.Range("L1").Value = "break_on_off"
$ws.Range("L1").Font.Bold = $true
$ws.Range("L1").HorizontalAlignment = -4108

# Rows where break_on_off = 1 (all others = 0)
$onRows = @(19, 37, 54)

for ($r = 2; $r -le 73; $r++) {
    if ($onRows -contains $r) {
        $ws.Cells.Item($r, 12).Value = 1
    } else {
        $ws.Cells.Item($r, 12).Value = 0
    }
}

# Update the selection / active cell to the new column, matching the edited file
$ws.Range("L1:L73").Select()
